# New run of the FIM (Fiscal Impact Model) - update computed columns L:S
# for rows 13, 16, 41, 44 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    13 = @{ L = 0.0044;  M = 0.0031;  N = 0.0041;  O = 0.0022;  P = 0.0022;  Q = 0.0018;  R = 0.0017;  S = -0.0035 }
    16 = @{ L = -0.3045; M = -0.2926; N = -0.5478; O = -0.0175; P = -0.0664; Q = 0.0577;  R = -0.0722; S = -72.7098 }
    41 = @{ L = -0.0393; M = -0.0277; N = -0.0359; O = -0.0195; P = -0.0192; Q = -0.016;  R = -0.0153; S = 0.0312 }
    44 = @{ L = -0.19;   M = -0.0595; N = -0.0897; O = 0.0765;  P = -0.057;  Q = -0.0487; R = -0.0418; S = -0.9983 }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$rowNum").Value = $cols[$col]
    }
}
